$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the current account (row 2, columns F:G) down to a new row 3,
# preserving its existing styling (hyperlink style on F, plain style on G).
$ws.Range("F2:G2").Copy($ws.Range("F3:G3"))

# Put the new account info into row 2.
$ws.Range("F2").Value = "yu.li9@hpe.com"
$ws.Range("G2").Value = "58c0df4d8413a6c2992ba35e7c56670f3ca3a56ed44bd79e4b5ae79ad0e5"

# Add a mailto hyperlink for the new account email, keeping the cell's
# original (Hyperlink) style instead of whatever the Add() call reapplies.
$origStyle = $ws.Range("F2").Style
$ws.Hyperlinks.Add($ws.Range("F2"), "mailto:yu.li9@hpe.com")
$ws.Range("F2").Style = $origStyle

# Match the saved selection/active cell from the source workbook.
$ws.Range("G9").Select()
